# GossA-HW10.xlsx: rename sheet, tighten a handful of quadrature results to
# higher precision, and append a new "HexGrid-60degTilt5degRes" row (row 16)
# produced by the new Gaussian Quadrature routine.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab to the short name.
$ws.Name = "GossA"

# --- tiny last-digit precision corrections from the re-run quadrature calc ---
$ws.Range("E13").Value = 1.005877367761607
$ws.Range("N13").Value = 0.9948687026895993
$ws.Range("O13").Value = 0.9918270117165819
$ws.Range("P13").Value = 0.992858039728209

$ws.Range("E15").Value = 0.9841436777953232
$ws.Range("I15").Value = 0.9644317021841152
$ws.Range("J15").Value = 0.9700759480451602

# --- new row 16: HexGrid-60degTilt5degRes averaged intensities ---
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 2.018785056086185
$ws.Range("D16").Value = 2.056121120826005
$ws.Range("E16").Value = 1.634508309208395
$ws.Range("F16").Value = 0.7388292315594466
$ws.Range("G16").Value = 2.018785056086185
$ws.Range("H16").Value = 2.056121120826005
$ws.Range("I16").Value = 1.050463488344989
$ws.Range("J16").Value = 0.5806850065064968
$ws.Range("K16").Value = 1.021648733139336
$ws.Range("L16").Value = 0.8781271138853387
$ws.Range("M16").Value = 2.018785056086185
$ws.Range("N16").Value = 1.8453147150172
$ws.Range("O16").Value = 1.612060929420008
$ws.Range("P16").Value = 1.247396007444524

# Column A on the data rows carries the bordered/bold/centered style (style
# index 1 in styles.xml) - copy it from the row above onto the new label cell.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$excel.CutCopyMode = 0
